$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 82, shifting rows 82:91 down to 83:92
$ws.Rows.Item(82).Insert()

# Fill in the new row 82 data
$ws.Cells.Item(82, 1).Value = 3
$ws.Cells.Item(82, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(82, 3).Value = "Coquimbo"
$ws.Cells.Item(82, 4).Value = 45124
$ws.Cells.Item(82, 5).Value = 5
$ws.Cells.Item(82, 6).Value = 100112022
$ws.Cells.Item(82, 7).Value = "Arveja Verde"
$ws.Cells.Item(82, 8).Value = "Perfection"
$ws.Cells.Item(82, 9).Value = "Primera"
$ws.Cells.Item(82, 10).Value = 73
$ws.Cells.Item(82, 11).Value = 28000
$ws.Cells.Item(82, 12).Value = 29000
$ws.Cells.Item(82, 13).Value = 28479
$ws.Cells.Item(82, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(82, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(82, 16).Value = 1139
$ws.Cells.Item(82, 17).Value = 25
$ws.Cells.Item(82, 18).Value = "Hortaliza"
